# New crime data collected — weekly 122nd Precinct CompStat update.
# Updates the report's issue/date header text and refreshes the crime-count
# table (rows 15-28) with the new weekly figures and recomputed percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: bump issue number (Volume 32 Number 8 -> 9) -------------------
$ws.Range("A8").Characters(21, 1).Text = "9"

# --- Header: shift the reporting week (2/17/2025-2/23/2025 -> 2/24/2025-3/2/2025) ---
$ws.Range("C9").Characters(27, 9).Text = "2/24/2025"
$ws.Range("C9").Characters(47, 9).Text = "3/2/2025"

# --- Row 15 (Rape) ----------------------------------------------------------
$ws.Range("G15").Value = 3
$ws.Range("J15").Value = 3

# --- Row 16 (Robbery) -------------------------------------------------------
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = -85.714285714285
$ws.Range("J16").Value = 13
$ws.Range("K16").Value = -53.846153846153
$ws.Range("L16").Value = -50
$ws.Range("N16").Value = -89.655172413793

# --- Row 17 (Fel. Assault) ---------------------------------------------------
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 60
$ws.Range("F17").Value = 12
$ws.Range("G17").Value = 18
$ws.Range("H17").Value = -33.333333333333
$ws.Range("I17").Value = 29
$ws.Range("J17").Value = 31
$ws.Range("K17").Value = -6.451612903225
$ws.Range("L17").Value = -6.451612903225
$ws.Range("M17").Value = 7.407407407407
$ws.Range("N17").Value = -32.558139534883

# --- Row 18 (Burglary) ------------------------------------------------------
$ws.Range("C18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = 28.571428571428
$ws.Range("I18").Value = 26
$ws.Range("J18").Value = 11
$ws.Range("K18").Value = 136.363636363636
$ws.Range("L18").Value = 62.5
$ws.Range("M18").Value = -18.75
$ws.Range("N18").Value = -88.235294117647

# --- Row 19 (Gr. Larceny) ---------------------------------------------------
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 26
$ws.Range("G19").Value = 38
$ws.Range("H19").Value = -31.578947368421
$ws.Range("I19").Value = 54
$ws.Range("J19").Value = 101
$ws.Range("K19").Value = -46.534653465346
$ws.Range("L19").Value = -39.325842696629
$ws.Range("M19").Value = -18.181818181818
$ws.Range("N19").Value = -61.151079136690

# --- Row 20 (G.L.A.) --------------------------------------------------------
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = -25
$ws.Range("I20").Value = 7
$ws.Range("J20").Value = 9
$ws.Range("K20").Value = -22.222222222222
$ws.Range("L20").Value = -61.111111111111
$ws.Range("M20").Value = -41.666666666666
$ws.Range("N20").Value = -98.501070663811

# --- Row 21 (TOTAL) ----------------------------------------------------------
$ws.Range("C21").Value = 15
$ws.Range("D21").Value = 16
$ws.Range("E21").Value = -6.25
$ws.Range("F21").Value = 51
$ws.Range("G21").Value = 77
$ws.Range("H21").Value = -33.766233766233
$ws.Range("I21").Value = 122
$ws.Range("J21").Value = 168
$ws.Range("K21").Value = -27.380952380952
$ws.Range("L21").Value = -26.506024096385
$ws.Range("M21").Value = -20.779220779220
$ws.Range("N21").Value = -86.895810955961

# --- Row 23 (Housing) -------------------------------------------------------
# C23 flips from the "0" placeholder text to a real numeric count; copy F23's
# number format onto it first so the cell's style matches the other numeric
# cells in the row, then write the value.
$ws.Range("F23").Copy() | Out-Null
$ws.Range("C23").PasteSpecial(-4122) | Out-Null
$ws.Range("C23").Value = 1
$ws.Range("F23").Value = 3
$ws.Range("I23").Value = 4
$ws.Range("L23").Value = -42.857142857142
$ws.Range("M23").Value = 33.333333333333

# --- Row 24 (Petit Larceny) -------------------------------------------------
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = -21.739130434782
$ws.Range("F24").Value = 87
$ws.Range("G24").Value = 97
$ws.Range("H24").Value = -10.309278350515
$ws.Range("I24").Value = 190
$ws.Range("J24").Value = 169
$ws.Range("K24").Value = 12.426035502958
$ws.Range("L24").Value = 17.283950617283
$ws.Range("M24").Value = -25.78125

# --- Row 25 (Retail Theft) --------------------------------------------------
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 15
$ws.Range("E25").Value = -40
$ws.Range("F25").Value = 46
$ws.Range("G25").Value = 57
$ws.Range("H25").Value = -19.298245614035
$ws.Range("I25").Value = 111
$ws.Range("J25").Value = 81
$ws.Range("K25").Value = 37.037037037037
$ws.Range("L25").Value = 81.967213114754

# --- Row 26 (Misd. Assault) -------------------------------------------------
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = -63.636363636363
$ws.Range("F26").Value = 11
$ws.Range("G26").Value = 33
$ws.Range("H26").Value = -66.666666666666
$ws.Range("I26").Value = 47
$ws.Range("J26").Value = 52
$ws.Range("K26").Value = -9.615384615384
$ws.Range("L26").Value = 34.285714285714
$ws.Range("M26").Value = -31.884057971014

# --- Row 27 (UCR Rape*) ------------------------------------------------------
$ws.Range("J27").Value = 4

# --- Row 28 (Other Sex Crimes) ----------------------------------------------
# C28 flips the other way: it was a numeric count and becomes the "0"
# placeholder text. Write the quoted text first, then copy D28's (already a
# "0" placeholder) format onto it so the cell's style matches.
$ws.Range("C28").Value = "'0"
$ws.Range("D28").Copy() | Out-Null
$ws.Range("C28").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 300
$ws.Range("L28").Value = 80
